$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("qw_e4")
Write-Host $ws.Name
